# Removed speaker notes from several slides.
$p = $ppt.ActivePresentation

$slideIndexes = @(5, 6, 7, 9, 10, 12, 13)

foreach ($idx in $slideIndexes) {
    $s = $p.Slides.Item($idx)
    $notesPage = $s.NotesPage

    $notesShape = $null
    for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
        $candidate = $notesPage.Shapes.Item($i)
        if ($candidate.Name -eq "Notes Placeholder 2") {
            $notesShape = $candidate
            break
        }
    }
    if ($notesShape -eq $null) {
        # Fallback: the notes text placeholder is conventionally shape 2.
        $notesShape = $notesPage.Shapes.Item(2)
    }

    $notesShape.TextFrame.TextRange.Text = ""
}
